$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1709.5385
$ws.Range("I70").Value = 1932.4
$ws.Range("J70").Value = 966.6667
$ws.Range("K70").Value = 5797.200000000001
$ws.Range("L70").Value = 2900.0001
$ws.Range("M70").Value = -5527.200000000001
$ws.Range("N70").Value = -3440.0001
$ws.Range("H73").Value = 1709.5385
$ws.Range("I73").Value = 1932.4
$ws.Range("J73").Value = 966.6667
$ws.Range("K73").Value = 5797.200000000001
$ws.Range("L73").Value = 2900.0001
$ws.Range("M73").Value = -4861.200000000001
$ws.Range("N73").Value = -4772.0001
$ws.Range("H98").Value = 54348944
$ws.Range("I98").Value = 78125940
$ws.Range("J98").Value = 1543.4286
$ws.Range("K98").Value = 78125940
$ws.Range("L98").Value = 1543.4286
$ws.Range("M98").Value = -78124442
$ws.Range("N98").Value = -4539.4286
$ws.Range("H100").Value = 10494.417
$ws.Range("I100").Value = 21645.4
$ws.Range("J100").Value = 2529.4285
$ws.Range("K100").Value = 21645.4
$ws.Range("L100").Value = 2529.4285
$ws.Range("M100").Value = -21104.4
$ws.Range("N100").Value = -3611.4285
$ws.Range("H110").Value = 62000
$ws.Range("J110").Value = 62000
$ws.Range("L110").Value = 62000
$ws.Range("N110").Value = -70180
$ws.Range("H113").Value = 3113.0435
$ws.Range("J113").Value = 2927.7778
$ws.Range("L113").Value = 2927.7778
$ws.Range("N113").Value = -9435.7778
$ws.Range("H116").Value = 2005
$ws.Range("I116").Value = 2005
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2005
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1437
$ws.Range("H117").Value = 80000
$ws.Range("J117").Value = 80000
$ws.Range("L117").Value = 80000
$ws.Range("N117").Value = -89178
$ws.Range("H122").Value = 54348944
$ws.Range("I122").Value = 78125940
$ws.Range("J122").Value = 1543.4286
$ws.Range("K122").Value = 234377820
$ws.Range("L122").Value = 4630.2858
$ws.Range("M122").Value = -234375370
$ws.Range("N122").Value = -9530.2858
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 84696.586
$ws.Range("I45").Value = 200772
$ws.Range("K45").Value = 200772
$ws.Range("M45").Value = -200395
$ws.Range("H61").Value = 2448.4443
$ws.Range("I61").Value = 1759
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1759
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1547
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 62502100
$ws.Range("I74").Value = 83335630
$ws.Range("J74").Value = 1499.5
$ws.Range("K74").Value = 83335630
$ws.Range("L74").Value = 1499.5
$ws.Range("M74").Value = -83334756
$ws.Range("N74").Value = -3247.5
$ws.Range("H77").Value = 62502100
$ws.Range("I77").Value = 83335630
$ws.Range("J77").Value = 1499.5
$ws.Range("K77").Value = 416678150
$ws.Range("L77").Value = 7497.5
$ws.Range("M77").Value = -416673782
$ws.Range("N77").Value = -16233.5
$ws.Range("H102").Value = 1305
$ws.Range("I102").Value = 1305
$ws.Range("K102").Value = 1305
$ws.Range("M102").Value = 317
$ws.Range("H110").Value = 1615.381
$ws.Range("I110").Value = 1042.8667
$ws.Range("K110").Value = 1042.8667
$ws.Range("M110").Value = 1002.1333
$ws.Range("H122").Value = 683.125
$ws.Range("I122").Value = 748.6842
$ws.Range("J122").Value = 434
$ws.Range("K122").Value = 2246.0526
$ws.Range("L122").Value = 1302
$ws.Range("M122").Value = 203.9474
$ws.Range("N122").Value = -6202
$ws.Range("H136").Value = 2448.4443
$ws.Range("I136").Value = 1759
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5277
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2727
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3122.5386
$ws.Range("I16").Value = 1333.3334
$ws.Range("J16").Value = 3659.3
$ws.Range("K16").Value = 1333.3334
$ws.Range("L16").Value = 3659.3
$ws.Range("M16").Value = -1046.3334
$ws.Range("N16").Value = -4233.3
$ws.Range("H31").Value = 1426.8733
$ws.Range("I31").Value = 1107.7778
$ws.Range("J31").Value = 1755.0857
$ws.Range("K31").Value = 1107.7778
$ws.Range("L31").Value = 1755.0857
$ws.Range("M31").Value = -812.7778000000001
$ws.Range("N31").Value = -2345.0857
$ws.Range("H34").Value = 1426.8733
$ws.Range("I34").Value = 1107.7778
$ws.Range("J34").Value = 1755.0857
$ws.Range("K34").Value = 1107.7778
$ws.Range("L34").Value = 1755.0857
$ws.Range("M34").Value = -905.7778000000001
$ws.Range("N34").Value = -2159.0857
$ws.Range("H99").Value = 142859630
$ws.Range("I99").Value = 333334880
$ws.Range("J99").Value = 3198.5
$ws.Range("K99").Value = 333334880
$ws.Range("L99").Value = 3198.5
$ws.Range("M99").Value = -333333382
$ws.Range("N99").Value = -6194.5
$ws.Range("H113").Value = 3122.5386
$ws.Range("I113").Value = 1333.3334
$ws.Range("J113").Value = 3659.3
$ws.Range("K113").Value = 1333.3334
$ws.Range("L113").Value = 3659.3
$ws.Range("M113").Value = 836.6666
$ws.Range("N113").Value = -7999.3
$ws.Range("H126").Value = 142859630
$ws.Range("I126").Value = 333334880
$ws.Range("J126").Value = 3198.5
$ws.Range("K126").Value = 1000004640
$ws.Range("L126").Value = 9595.5
$ws.Range("M126").Value = -1000002170
$ws.Range("N126").Value = -14535.5
$ws.Range("H132").Value = 10418475
$ws.Range("I132").Value = 1560.5714
$ws.Range("K132").Value = 4681.7142
$ws.Range("M132").Value = -2151.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 6677545
$ws.Range("I22").Value = 50000750
$ws.Range("J22").Value = 12436.462
$ws.Range("K22").Value = 150002250
$ws.Range("L22").Value = 37309.386
$ws.Range("M22").Value = -150002081
$ws.Range("N22").Value = -37647.386
$ws.Range("H27").Value = 6677545
$ws.Range("I27").Value = 50000750
$ws.Range("J27").Value = 12436.462
$ws.Range("K27").Value = 150002250
$ws.Range("L27").Value = 37309.386
$ws.Range("M27").Value = -150002148
$ws.Range("N27").Value = -37513.386
$ws.Range("H33").Value = 4399.56
$ws.Range("I33").Value = 815.8889
$ws.Range("J33").Value = 6415.375
$ws.Range("K33").Value = 4895.3334
$ws.Range("L33").Value = 38492.25
$ws.Range("M33").Value = -4612.3334
$ws.Range("N33").Value = -39058.25
$ws.Range("H34").Value = 713
$ws.Range("I34").Value = 450.66666
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1351.99998
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -1267.99998
$ws.Range("N34").Value = -4668
$ws.Range("H131").Value = 780.4897999999999
$ws.Range("J131").Value = 797.79346
$ws.Range("L131").Value = 2393.38038
$ws.Range("N131").Value = -12473.38038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4383.12
$ws.Range("I107").Value = 253.63637
$ws.Range("J107").Value = 7627.7144
$ws.Range("K107").Value = 253.63637
$ws.Range("L107").Value = 7627.7144
$ws.Range("M107").Value = 1666.36363
$ws.Range("N107").Value = -11467.7144
$ws.Range("H113").Value = 877.8
$ws.Range("I113").Value = 825.4286
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 825.4286
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1344.5714
$ws.Range("N113").Value = -5340
$ws.Range("H122").Value = 35722100
$ws.Range("I122").Value = 35722100
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 107166300
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -107163850
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1637.6842
$ws.Range("I7").Value = 1579.7142
$ws.Range("K7").Value = 1579.7142
$ws.Range("M7").Value = -1467.7142
$ws.Range("H40").Value = 31251808
$ws.Range("I40").Value = 2100
$ws.Range("J40").Value = 83334660
$ws.Range("K40").Value = 2100
$ws.Range("L40").Value = 83334660
$ws.Range("M40").Value = -1964
$ws.Range("N40").Value = -83334932
$ws.Range("H126").Value = 1637.6842
$ws.Range("I126").Value = 1579.7142
$ws.Range("K126").Value = 4739.142599999999
$ws.Range("M126").Value = -2269.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1491.7142
$ws.Range("I100").Value = 1212.5834
$ws.Range("J100").Value = 3166.5
$ws.Range("K100").Value = 2425.1668
$ws.Range("L100").Value = 6333
$ws.Range("M100").Value = -1884.1668
$ws.Range("N100").Value = -7415
$ws.Range("H126").Value = 1490.2667
$ws.Range("I126").Value = 999
$ws.Range("J126").Value = 1817.7778
$ws.Range("K126").Value = 2997
$ws.Range("L126").Value = 5453.3334
$ws.Range("M126").Value = -527
$ws.Range("N126").Value = -10393.3334
